$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.026.15"
$ws.Range("E2").Value = "  +2.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.80"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.27"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9968"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3925"
$ws.Range("E7").Value = "  +3.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3486"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.47"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.203"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07598"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9965"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.21"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.539"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.812.50"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001109"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06704"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.22"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9974"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("E21").Value = "  +2.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.577"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.017.61"
$ws.Range("E23").Value = "  +2.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.84"
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.401"
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.526"
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.574"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.39"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "155.00"
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.038.63"
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.25"
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.041"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.152"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08850"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.33"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.534"
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6964"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02428"
$ws.Range("E38").Value = "  +3.68%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06543"
$ws.Range("E39").Value = "  +2.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.613"
$ws.Range("E40").Value = "  -4.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2237"
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.267"
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.563"
$ws.Range("E43").Value = "  -3.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.74"
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6535"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9963"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.874"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.174"
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.82"
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07223"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.42"
$ws.Range("E51").Value = "  +1.37%  "
